# Add inferred sequence documentation
#
# 1. Remove the duplicate "Question for review..." note from the
#    Inferences sheet (it already appears, worded slightly differently,
#    on the Genotype sheet).
# 2. Add a new "Acknowledgements" explanation block to the Submission
#    sheet, mirroring the existing sections on the other tabs.
# 3. Leave the cursor/active tab on the Submission sheet, matching the
#    selections left behind by the edit.

$wb = $excel.ActiveWorkbook

# --- 1. Inferences sheet: drop the redundant reviewer question row ---
$wsInferences = $wb.Worksheets.Item("Inferences")
$wsInferences.Range("B19").EntireRow.Delete()

# --- 2. Submission sheet: new "Acknowledgements" notes section ---
$wsSubmission = $wb.Worksheets.Item("Submission")

$wsSubmission.Range("B16").Value = "Acknowledgements"
$wsSubmission.Range("B16").Font.Bold = $true

$wsSubmission.Range("B17").Value = "Please list the individuals whose contribution to this work should be acknowledged"

$wsSubmission.Range("B19").Value = "{{Acknowledgements:properties#10}}"

# Write B22 before B21 so the new shared strings are appended to the
# table in the same order as the source workbook.
$wsSubmission.Range("B22").Value = "{{Acknowledgements:properties!type,completed_by}}"

$wsSubmission.Range("B21").Value = "Explanation of fields"
$wsSubmission.Range("B21").Font.Italic = $true

# --- 3. Selections / active sheet left by the edit ---
$wsRepertoire = $wb.Worksheets.Item("Repertoire")
$wsRepertoire.Range("B11").Select()

$wsInferences.Range("B21").Select()

$wsSubmission.Activate()
$wsSubmission.Range("B21").Select()
